# modals for volatile data acquisition
# Insert 4 new key/value rows ("status", "cpu", "memory", "network") before the
# existing "malware_check_title" row (old row 54), and append one new row
# ("cache_data_captured") at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four blank rows at row 54, pushing the existing data (old rows 54-77)
# down to rows 58-81.
$ws.Rows("54:57").Insert()

# Fill in the four newly-inserted rows.
$ws.Range("A54").Value = "status"
$ws.Range("B54").Value = "Status"

$ws.Range("A55").Value = "cpu"
$ws.Range("B55").Value = "CPU"

$ws.Range("A56").Value = "memory"
$ws.Range("B56").Value = "Memory"

$ws.Range("A57").Value = "network"
$ws.Range("B57").Value = "Network"

# Append a new row (82) at the bottom of the sheet.
$ws.Range("A82").Value = "cache_data_captured"
$ws.Range("B82").Value = "Cache data copied: command history, clipboard, print spool files."

# Match the author's final selection/view position.
$null = $ws.Range("A82").Select()
